{"js": "// Replace the date line and each division-problem cell text with its\n// updated value. Every old string is unique within the document, so a\n// simple search + replace (matching the whole text, case-sensitive) is\n// sufficient and safe.\nconst replacements = [\n  [\"2025-06-26 Thursday\", \"2025-06-27 Friday\"],\n  [\"27\u00f76=4, 3\", \"12\u00f79=1, 3\"],\n  [\"33\u00f77=4, 5\", \"62\u00f75=12, 2\"],\n  [\"49\u00f77=7, 0\", \"86\u00f79=9, 5\"],\n  [\"88\u00f77=12, 4\", \"28\u00f75=5, 3\"],\n  [\"86\u00f74=21, 2\", \"41\u00f72=20, 1\"],\n  [\"75\u00f73=25, 0\", \"37\u00f79=4, 1\"],\n  [\"34\u00f74=8, 2\", \"39\u00f78=4, 7\"],\n  [\"10\u00f74=2, 2\", \"59\u00f76=9, 5\"],\n  [\"22\u00f73=7, 1\", \"54\u00f77=7, 5\"],\n  [\"47\u00f76=7, 5\", \"89\u00f72=44, 1\"],\n  [\"32\u00f76=5, 2\", \"40\u00f73=13, 1\"],\n  [\"30\u00f72=15, 0\", \"85\u00f73=28, 1\"],\n  [\"92\u00f72=46, 0\", \"22\u00f75=4, 2\"],\n  [\"75\u00f79=8, 3\", \"68\u00f72=34, 0\"],\n  [\"37\u00f78=4, 5\", \"92\u00f75=18, 2\"],\n  [\"60\u00f75=12, 0\", \"73\u00f75=14, 3\"],\n  [\"48\u00f79=5, 3\", \"84\u00f77=12, 0\"],\n  [\"65\u00f72=32, 1\", \"61\u00f73=20, 1\"],\n  [\"99\u00f77=14, 1\", \"42\u00f76=7, 0\"],\n  [\"16\u00f77=2, 2\", \"39\u00f73=13, 0\"],\n  [\"17\u00f78=2, 1\", \"59\u00f79=6, 5\"],\n  [\"70\u00f79=7, 7\", \"72\u00f75=14, 2\"],\n  [\"91\u00f72=45, 1\", \"53\u00f77=7, 4\"],\n  [\"18\u00f79=2, 0\", \"74\u00f77=10, 4\"],\n  [\"20\u00f74=5, 0\", \"94\u00f79=10, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each \"find\" string is unique in the document (the date line plus one\n# string per table cell containing a division problem), so a simple\n# Find/Replace (ReplaceAll) scoped to the whole document body is safe\n# and replaces exactly one occurrence each.\n$replacements = @(\n    @(\"2025-06-26 Thursday\", \"2025-06-27 Friday\"),\n    @(\"27\u00f76=4, 3\", \"12\u00f79=1, 3\"),\n    @(\"33\u00f77=4, 5\", \"62\u00f75=12, 2\"),\n    @(\"49\u00f77=7, 0\", \"86\u00f79=9, 5\"),\n    @(\"88\u00f77=12, 4\", \"28\u00f75=5, 3\"),\n    @(\"86\u00f74=21, 2\", \"41\u00f72=20, 1\"),\n    @(\"75\u00f73=25, 0\", \"37\u00f79=4, 1\"),\n    @(\"34\u00f74=8, 2\", \"39\u00f78=4, 7\"),\n    @(\"10\u00f74=2, 2\", \"59\u00f76=9, 5\"),\n    @(\"22\u00f73=7, 1\", \"54\u00f77=7, 5\"),\n    @(\"47\u00f76=7, 5\", \"89\u00f72=44, 1\"),\n    @(\"32\u00f76=5, 2\", \"40\u00f73=13, 1\"),\n    @(\"30\u00f72=15, 0\", \"85\u00f73=28, 1\"),\n    @(\"92\u00f72=46, 0\", \"22\u00f75=4, 2\"),\n    @(\"75\u00f79=8, 3\", \"68\u00f72=34, 0\"),\n    @(\"37\u00f78=4, 5\", \"92\u00f75=18, 2\"),\n    @(\"60\u00f75=12, 0\", \"73\u00f75=14, 3\"),\n    @(\"48\u00f79=5, 3\", \"84\u00f77=12, 0\"),\n    @(\"65\u00f72=32, 1\", \"61\u00f73=20, 1\"),\n    @(\"99\u00f77=14, 1\", \"42\u00f76=7, 0\"),\n    @(\"16\u00f77=2, 2\", \"39\u00f73=13, 0\"),\n    @(\"17\u00f78=2, 1\", \"59\u00f79=6, 5\"),\n    @(\"70\u00f79=7, 7\", \"72\u00f75=14, 2\"),\n    @(\"91\u00f72=45, 1\", \"53\u00f77=7, 4\"),\n    @(\"18\u00f79=2, 0\", \"74\u00f77=10, 4\"),\n    @(\"20\u00f74=5, 0\", \"94\u00f79=10, 4\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n$d.Save()\n"}
